# Update "想去人数" (interest counts) for two events that were refreshed
# in the source data, mirroring the values across both the "展览"
# (exhibition) sheet and the "全部类型" (all types) sheet.

$wb = $excel.ActiveWorkbook

# "展览" sheet: rows 4 and 5 hold the two affected events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 4912
$wsExhibit.Range("F5").Value = 25

# "全部类型" sheet: same two events, but located at rows 4 and 6
# because of an extra row inserted between them on this sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 4912
$wsAll.Range("F6").Value = 25
